$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New translated-text strings (added to xl/sharedStrings.xml as indices
# 71..89). They must be assigned to cells in this exact order so that the
# shared-string table receives new entries in the same sequence as the
# target workbook.
# ---------------------------------------------------------------------------
${s71} = ' I had a Pulpy Life Seed earlier,\nand it was pretty tasty.'
${s72} = ' SCRIPT/P01P04A/us0404.ssb'
${s73} = ' Недавно я пил Сок из Зерна Жизни.\nОн довольно вкусный.'
${s74} = ' Îåäàâîï ÿ ðéì Òïë éè Èåñîà Çéèîé.\nÏî äïâïìûîï âëôòîúê.'
${s75} = 'SCRIPT/P01P04A/us2009.ssb'
${s76} = ' We heard that [CS:N]Drowzee[CR] is holed\nup on [CS:P]Mt. Travail[CR].'
${s77} = ' Мы слышали, что [CS:N]Дроузи[CR] окопался\nна [CS:P]Горе Травейл[CR].'
${s78} = 'SCRIPT/T01P02A/us2011.ssb '
${s79} = ' Íú òìúšàìé, œóï [CS:N]Äñïôèé[CR] ïëïðàìòÿ\nîà [CS:P]Ãïñå Óñàâåêì[CR].'
${s80} = ' [CS:N]Drowzee[CR] has become a changed\nPokémon?[K] Good to hear he has reformed.'
${s81} = ' Not that it matters to us, since\nwe''re a team that focuses on catching outlaws.'
${s82} = ' But if [CS:N]Drowzee[CR] becomes an\noutlaw again, we''ll be there to catch him.'
${s83} = 'SCRIPT/T01P02A/us2015.ssb'
${s84} = ' Не то чтобы это нас заботит,\nведь мы команда, которая ловит негодяев.'
${s85} = ' Но если [CS:N]Дроузи[CR] снова станет\nнегодяем, мы обязательно поймаем его.'
${s86} = ' Îå óï œóïáú üóï îàò èàáïóéó,\nâåäû íú ëïíàîäà, ëïóïñàÿ ìïâéó îåãïäÿåâ.'
${s87} = ' Îï åòìé [CS:N]Äñïôèé[CR] òîïâà òóàîåó\nîåãïäÿåí, íú ïáÿèàóåìûîï ðïêíàåí åãï.'
${s88} = ' [CS:N]Äñïôèé[CR] éèíåîéìòÿ?[K] Ðñéÿóîï èîàóû,\nœóï ïî éòðñàâéìòÿ.'
${s89} = ' [CS:N]Дроузи[CR] изменился?[K] Приятно знать,\nчто он исправился.'

# ---------------------------------------------------------------------------
# Row 21 changes its formatting group (style 4/5 -> 6/7) and gains an empty,
# styled A21 cell. Copy the formats only (values are left untouched) from an
# existing row that already uses the target style group.
# ---------------------------------------------------------------------------
$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Add six new rows (22-27) of translated script content below the existing
# table, copying cell formatting from rows that already carry the right
# style group before filling in the values.
# ---------------------------------------------------------------------------
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null

$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A23:E23").PasteSpecial(-4122) | Out-Null

$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A24:E24").PasteSpecial(-4122) | Out-Null

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A25:E25").PasteSpecial(-4122) | Out-Null

$ws.Range("B3:E3").Copy() | Out-Null
$ws.Range("B26:E26").PasteSpecial(-4122) | Out-Null

$ws.Range("B3:E3").Copy() | Out-Null
$ws.Range("B27:E27").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Fill in the cell values, in the exact order that produces shared-string
# indices 71..89.
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = ${s71}
$ws.Range("A22").Value = ${s72}
$ws.Range("D22").Value = ${s73}
$ws.Range("E22").Value = ${s74}

$ws.Range("A23").Value = ${s75}

$ws.Range("C24").Value = ${s76}
$ws.Range("D24").Value = ${s77}
$ws.Range("A24").Value = ${s78}
$ws.Range("E24").Value = ${s79}

$ws.Range("C25").Value = ${s80}
$ws.Range("C26").Value = ${s81}
$ws.Range("C27").Value = ${s82}

$ws.Range("A25").Value = ${s83}

$ws.Range("D26").Value = ${s84}
$ws.Range("D27").Value = ${s85}
$ws.Range("E26").Value = ${s86}
$ws.Range("E27").Value = ${s87}
$ws.Range("E25").Value = ${s88}
$ws.Range("D25").Value = ${s89}

# Numeric "row number" column values.
$ws.Range("B22").Value = 63
$ws.Range("B24").Value = 44
$ws.Range("B25").Value = 18
$ws.Range("B26").Value = 21
$ws.Range("B27").Value = 24

# ---------------------------------------------------------------------------
# Row heights for the newly added rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 57.6
$ws.Rows.Item(23).RowHeight = 43.2
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 43.2
$ws.Rows.Item(26).RowHeight = 21.6
$ws.Rows.Item(27).RowHeight = 31.8

# ---------------------------------------------------------------------------
# Update the view: scroll so row 22 is at the top and select E25, matching
# the editor's on-screen state when the change was committed.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E25").Select() | Out-Null
